$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column A (key names) with a narrower width; column B keeps its existing width
$ws.Columns.Item(1).ColumnWidth = 19.46

$ws.Range("A3").Value = 'start_controls1'
$ws.Range("B3").Value = 'Sometimes games with the simplest control schemes are the most difficult.  For example, in Getting Over It with Bennett Foddy, the only thing you do is move the mouse, but the complex interactions with physics allow for very in-depth movement mechanics. '

$ws.Range("A4").Value = 'start_controls2'
$ws.Range("B4").Value = 'In this game, designed to capture the Foddian spirit, you simply have 2 controls: flop left and flop right.  You can use arrow keys, A and D, or the mouse buttons to flip and flop.'

$ws.Range("A5").Value = 'start_controls3'
$ws.Range("B5").Value = 'You are, quite literally, a fish out of water.  I was planning to have the the movement driven entirely by the physics of the fish contacting the world, but, believe it or not, that was more uncontrollable than what you’re playing now.  You’d flip and flop upward with no practical way to progress forward.'

$ws.Range("A6").Value = 'start_controls4'
$ws.Range("B6").Value = 'Then I realized that, while not as effective as in water, a tail flapping back and forth would create some sort of forward momentum due to the air being pushed.'

$ws.Range("A7").Value = 'start_controls5'
$ws.Range("B7").Value = 'I’m not sure how physically accurate this is, but it does seem slightly more plausible than a man lifting his entire body weight and that of a cauldron up using just a hammer.'

$ws.Range("B6").Select()
